# ---------------------------------------------------------------------------
# "completed new opportunity scripts"
#
# 1. GeneralConfig: bump the iOS device/version parameters used by the
#    Appium test config (deviceNameIOS / platformVersionIOS).
# 2. ListingInfoTest: a new verification step (verifyFillingLiveWeightDetails)
#    is inserted right before the final "verifyPublishListing" step, pushing
#    that row down.
# 3. A brand-new worksheet (ListingInfoClassifiedTest) is added at the end,
#    mirroring the layout of ListingInfoTest, to drive the new "Classified"
#    listing flow test cases.
# 4. Misc. sheet selections / the active tab are left where the author's
#    Excel session ended up (LoginTest).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. GeneralConfig: iOS device params ----------------------------------
$general = $wb.Worksheets.Item("GeneralConfig")
$general.Range("B8").Value  = "iPad Air (3rd generation)"
$general.Range("B9").Value  = 13.4

# --- 2. ListingInfoTest: insert the new "live weight" verification row ----
$listing = $wb.Worksheets.Item("ListingInfoTest")

# Duplicate the formatting of the last existing row (13 - verifyPublishListing)
# down into the new row 14 so the new row inherits the same cell styles.
$listing.Range("A13:C13").Copy()
$listing.Range("A14:C14").PasteSpecial(-4122)
$listing.Rows.Item(14).RowHeight = 22

# Row 13 becomes the new verification step; the old "verifyPublishListing"
# step (and its "run" param) moves down to row 14.
$listing.Cells.Item(13, 1).Value = "verifyFillingLiveWeightDetails"
$listing.Cells.Item(14, 1).Value = "verifyPublishListing"
$listing.Cells.Item(14, 2).Value = "run"

# --- 3. New sheet: ListingInfoClassifiedTest -------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$classified = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$classified.Name = "ListingInfoClassifiedTest"

# Mirror column widths / header+row formatting from ListingInfoTest.
$classified.Columns.Item(1).ColumnWidth = 65.14
$classified.Columns.Item(2).ColumnWidth = 18
$classified.Columns.Item(3).ColumnWidth = 70.29

$listing.Range("A1:C1").Copy()
$classified.Range("A1:C1").PasteSpecial(-4122)
$classified.Rows.Item(1).RowHeight = 22

$listing.Range("A2:C2").Copy()
$classified.Range("A2:C8").PasteSpecial(-4122)
$classified.Rows.Item("2:8").RowHeight = 22

$classified.Cells.Item(1, 1).Value = "Test Case Name"
$classified.Cells.Item(1, 2).Value = "Run Status"
$classified.Cells.Item(1, 3).Value = "Params"

$classifiedSteps = @(
    "CreateOpportunityForClassified",
    "verifyListingTypeSelectionClassified",
    "verifyListingInfoNavigationClassified",
    "verifyFillingClassifiedListingOverviewDetails",
    "verifyFillingClassifiedLotDetails",
    "verifyFillingClassifiedHealthVetDetails",
    "verifyClassifiedPublishListing"
)
for ($i = 0; $i -lt $classifiedSteps.Count; $i++) {
    $row = $i + 2
    $classified.Cells.Item($row, 1).Value = $classifiedSteps[$i]
    $classified.Cells.Item($row, 2).Value = "run"
}
$classified.Range("D8").Select()

# --- 4. Leave the sheet selections / active tab the way the author did ----
$general.Range("B23").Select()

$login = $wb.Worksheets.Item("LoginTest")
$login.Range("A28").Select()

$newOpportunity = $wb.Worksheets.Item("NewOpportunityTest")
$newOpportunity.Range("A24").Select()

$dashboard = $wb.Worksheets.Item("DashboardTest")
$dashboard.Range("B21").Select()

$listing.Range("B25").Select()

# LoginTest is the tab left active.
$login.Activate()
